# Auto-generated Excel COM-interop script
$wb = $excel.ActiveWorkbook

$level0 = $wb.Worksheets.Item("Level0")
    $lvl0_row4 = New-Object 'object[,]' 1,9
    $lvl0_row4[0,0] = 990
    $lvl0_row4[0,1] = 1147
    $lvl0_row4[0,2] = 918
    $lvl0_row4[0,3] = 1017
    $lvl0_row4[0,4] = 1075
    $lvl0_row4[0,5] = 944
    $lvl0_row4[0,6] = 3055
    $lvl0_row4[0,7] = 3036
    $lvl0_row4[0,8] = 6091
$level0.Range("B4:J4").Value = $lvl0_row4

    $lvl0_row6 = New-Object 'object[,]' 1,9
    $lvl0_row6[0,0] = 770
    $lvl0_row6[0,1] = 550
    $lvl0_row6[0,2] = 673
    $lvl0_row6[0,3] = 612
    $lvl0_row6[0,4] = 784
    $lvl0_row6[0,5] = 550
    $lvl0_row6[0,6] = 1993
    $lvl0_row6[0,7] = 1946
    $lvl0_row6[0,8] = 3939
$level0.Range("B6:J6").Value = $lvl0_row6

$level0prop = $wb.Worksheets.Add($null, $level0)
$level0prop.Name = "Level0-prop"

$level0prop.Range("B1:G1").NumberFormat = "@"
    $lvl0_prop_data = New-Object 'object[,]' 6,10
    $lvl0_prop_data[0,0] = "Cell Type"
    $lvl0_prop_data[0,1] = "2262"
    $lvl0_prop_data[0,2] = "2263"
    $lvl0_prop_data[0,3] = "2264"
    $lvl0_prop_data[0,4] = "2265"
    $lvl0_prop_data[0,5] = "2266"
    $lvl0_prop_data[0,6] = "2267"
    $lvl0_prop_data[0,7] = "CTRL"
    $lvl0_prop_data[0,8] = "MM"
    $lvl0_prop_data[0,9] = "total"
    $lvl0_prop_data[1,0] = "Endothelium"
    $lvl0_prop_data[1,1] = 0.008282716731087797
    $lvl0_prop_data[1,2] = 0.007479861910241657
    $lvl0_prop_data[1,3] = 0.015188335358444714
    $lvl0_prop_data[1,4] = 0.010752688172043012
    $lvl0_prop_data[1,5] = 0.005780346820809248
    $lvl0_prop_data[1,6] = 0.024936061381074168
    $lvl0_prop_data[1,7] = 0.010202117420596728
    $lvl0_prop_data[1,8] = 0.013226998638397198
    $lvl0_prop_data[1,9] = 0.011706656346749226
    $lvl0_prop_data[2,0] = "Epithelium"
    $lvl0_prop_data[2,1] = 0.00717835450027609
    $lvl0_prop_data[2,2] = 0.00805523590333717
    $lvl0_prop_data[2,3] = 0.010935601458080195
    $lvl0_prop_data[2,4] = 0.005973715651135006
    $lvl0_prop_data[2,5] = 0.010509721492380452
    $lvl0_prop_data[2,6] = 0.008951406649616368
    $lvl0_prop_data[2,7] = 0.008662175168431183
    $lvl0_prop_data[2,8] = 0.008558646177786422
    $lvl0_prop_data[2,9] = 0.008610681114551084
    $lvl0_prop_data[3,0] = "Granulosa"
    $lvl0_prop_data[3,1] = 0.5466593042517945
    $lvl0_prop_data[3,2] = 0.6599539700805523
    $lvl0_prop_data[3,3] = 0.5577156743620899
    $lvl0_prop_data[3,4] = 0.6075268817204301
    $lvl0_prop_data[3,5] = 0.5648975302154493
    $lvl0_prop_data[3,6] = 0.6035805626598465
    $lvl0_prop_data[3,7] = 0.588065447545717
    $lvl0_prop_data[3,8] = 0.5905465862672632
    $lvl0_prop_data[3,9] = 0.5892995356037152
    $lvl0_prop_data[4,0] = "Immune"
    $lvl0_prop_data[4,1] = 0.012700165654334622
    $lvl0_prop_data[4,2] = 0.00805523590333717
    $lvl0_prop_data[4,3] = 0.007290400972053463
    $lvl0_prop_data[4,4] = 0.01015531660692951
    $lvl0_prop_data[4,5] = 0.006831318970047294
    $lvl0_prop_data[4,6] = 0.010869565217391304
    $lvl0_prop_data[4,7] = 0.009432146294513956
    $lvl0_prop_data[4,8] = 0.00914219023536277
    $lvl0_prop_data[4,9] = 0.009287925696594427
    $lvl0_prop_data[5,0] = "Mesenchyme"
    $lvl0_prop_data[5,1] = 0.4251794588625069
    $lvl0_prop_data[5,2] = 0.31645569620253167
    $lvl0_prop_data[5,3] = 0.4088699878493317
    $lvl0_prop_data[5,4] = 0.3655913978494624
    $lvl0_prop_data[5,5] = 0.4119810825013137
    $lvl0_prop_data[5,6] = 0.3516624040920716
    $lvl0_prop_data[5,7] = 0.3836381135707411
    $lvl0_prop_data[5,8] = 0.3785255786811904
    $lvl0_prop_data[5,9] = 0.3810952012383901
$level0prop.Range("A1:J6").Value = $lvl0_prop_data

$level0prop.Range("A1:J1").Font.Bold = $true
$level0prop.Range("A1:J1").HorizontalAlignment = -4108

$level1 = $wb.Worksheets.Item("Level1")
    $lvl1_row4 = New-Object 'object[,]' 1,9
    $lvl1_row4[0,0] = 50
    $lvl1_row4[0,1] = 38
    $lvl1_row4[0,2] = 43
    $lvl1_row4[0,3] = 23
    $lvl1_row4[0,4] = 30
    $lvl1_row4[0,5] = 37
    $lvl1_row4[0,6] = 131
    $lvl1_row4[0,7] = 90
    $lvl1_row4[0,8] = 221
$level1.Range("B4:J4").Value = $lvl1_row4

    $lvl1_row5 = New-Object 'object[,]' 1,9
    $lvl1_row5[0,0] = 52
    $lvl1_row5[0,1] = 25
    $lvl1_row5[0,2] = 46
    $lvl1_row5[0,3] = 47
    $lvl1_row5[0,4] = 56
    $lvl1_row5[0,5] = 36
    $lvl1_row5[0,6] = 123
    $lvl1_row5[0,7] = 139
    $lvl1_row5[0,8] = 262
$level1.Range("B5:J5").Value = $lvl1_row5

    $lvl1_row6 = New-Object 'object[,]' 1,9
    $lvl1_row6[0,0] = 353
    $lvl1_row6[0,1] = 219
    $lvl1_row6[0,2] = 207
    $lvl1_row6[0,3] = 295
    $lvl1_row6[0,4] = 146
    $lvl1_row6[0,5] = 157
    $lvl1_row6[0,6] = 779
    $lvl1_row6[0,7] = 598
    $lvl1_row6[0,8] = 1377
$level1.Range("B6:J6").Value = $lvl1_row6

    $lvl1_row7 = New-Object 'object[,]' 1,9
    $lvl1_row7[0,0] = 124
    $lvl1_row7[0,1] = 228
    $lvl1_row7[0,2] = 135
    $lvl1_row7[0,3] = 225
    $lvl1_row7[0,4] = 289
    $lvl1_row7[0,5] = 207
    $lvl1_row7[0,6] = 487
    $lvl1_row7[0,7] = 721
    $lvl1_row7[0,8] = 1208
$level1.Range("B7:J7").Value = $lvl1_row7

    $lvl1_row8 = New-Object 'object[,]' 1,9
    $lvl1_row8[0,0] = 411
    $lvl1_row8[0,1] = 637
    $lvl1_row8[0,2] = 487
    $lvl1_row8[0,3] = 427
    $lvl1_row8[0,4] = 554
    $lvl1_row8[0,5] = 507
    $lvl1_row8[0,6] = 1535
    $lvl1_row8[0,7] = 1488
    $lvl1_row8[0,8] = 3023
$level1.Range("B8:J8").Value = $lvl1_row8

    $lvl1_row11 = New-Object 'object[,]' 1,9
    $lvl1_row11[0,0] = 101
    $lvl1_row11[0,1] = 129
    $lvl1_row11[0,2] = 102
    $lvl1_row11[0,3] = 110
    $lvl1_row11[0,4] = 205
    $lvl1_row11[0,5] = 100
    $lvl1_row11[0,6] = 332
    $lvl1_row11[0,7] = 415
    $lvl1_row11[0,8] = 747
$level1.Range("B11:J11").Value = $lvl1_row11

    $lvl1_row12 = New-Object 'object[,]' 1,9
    $lvl1_row12[0,0] = 15
    $lvl1_row12[0,1] = 9
    $lvl1_row12[0,2] = 7
    $lvl1_row12[0,3] = 13
    $lvl1_row12[0,4] = 9
    $lvl1_row12[0,5] = 6
    $lvl1_row12[0,6] = 31
    $lvl1_row12[0,7] = 28
    $lvl1_row12[0,8] = 59
$level1.Range("B12:J12").Value = $lvl1_row12

    $lvl1_row15 = New-Object 'object[,]' 1,9
    $lvl1_row15[0,0] = 377
    $lvl1_row15[0,1] = 215
    $lvl1_row15[0,2] = 299
    $lvl1_row15[0,3] = 231
    $lvl1_row15[0,4] = 319
    $lvl1_row15[0,5] = 210
    $lvl1_row15[0,6] = 891
    $lvl1_row15[0,7] = 760
    $lvl1_row15[0,8] = 1651
$level1.Range("B15:J15").Value = $lvl1_row15

    $lvl1_row16 = New-Object 'object[,]' 1,9
    $lvl1_row16[0,0] = 225
    $lvl1_row16[0,1] = 154
    $lvl1_row16[0,2] = 192
    $lvl1_row16[0,3] = 228
    $lvl1_row16[0,4] = 204
    $lvl1_row16[0,5] = 189
    $lvl1_row16[0,6] = 571
    $lvl1_row16[0,7] = 621
    $lvl1_row16[0,8] = 1192
$level1.Range("B16:J16").Value = $lvl1_row16

$level1prop = $wb.Worksheets.Add($null, $level1)
$level1prop.Name = "Level1-prop"

$level1prop.Range("B1:G1").NumberFormat = "@"
    $lvl1_prop_data = New-Object 'object[,]' 18,10
    $lvl1_prop_data[0,0] = "Cell Type"
    $lvl1_prop_data[0,1] = "2262"
    $lvl1_prop_data[0,2] = "2263"
    $lvl1_prop_data[0,3] = "2264"
    $lvl1_prop_data[0,4] = "2265"
    $lvl1_prop_data[0,5] = "2266"
    $lvl1_prop_data[0,6] = "2267"
    $lvl1_prop_data[0,7] = "CTRL"
    $lvl1_prop_data[0,8] = "MM"
    $lvl1_prop_data[0,9] = "total"
    $lvl1_prop_data[1,0] = "Endothelium"
    $lvl1_prop_data[1,1] = 0.12096774193548387
    $lvl1_prop_data[1,2] = 0.05701754385964912
    $lvl1_prop_data[1,3] = 0.18518518518518517
    $lvl1_prop_data[1,4] = 0.08
    $lvl1_prop_data[1,5] = 0.03806228373702422
    $lvl1_prop_data[1,6] = 0.18840579710144928
    $lvl1_prop_data[1,7] = 0.10882956878850103
    $lvl1_prop_data[1,8] = 0.09431345353675451
    $lvl1_prop_data[1,9] = 0.10016556291390728
    $lvl1_prop_data[2,0] = "Epithelium"
    $lvl1_prop_data[2,1] = 0.10483870967741936
    $lvl1_prop_data[2,2] = 0.06140350877192982
    $lvl1_prop_data[2,3] = 0.13333333333333333
    $lvl1_prop_data[2,4] = 0.044444444444444446
    $lvl1_prop_data[2,5] = 0.06920415224913495
    $lvl1_prop_data[2,6] = 0.06763285024154589
    $lvl1_prop_data[2,7] = 0.09240246406570841
    $lvl1_prop_data[2,8] = 0.06102635228848821
    $lvl1_prop_data[2,9] = 0.07367549668874172
    $lvl1_prop_data[3,0] = "GC_Atretic"
    $lvl1_prop_data[3,1] = 0.4032258064516129
    $lvl1_prop_data[3,2] = 0.16666666666666666
    $lvl1_prop_data[3,3] = 0.31851851851851853
    $lvl1_prop_data[3,4] = 0.10222222222222223
    $lvl1_prop_data[3,5] = 0.10380622837370242
    $lvl1_prop_data[3,6] = 0.178743961352657
    $lvl1_prop_data[3,7] = 0.26899383983572894
    $lvl1_prop_data[3,8] = 0.12482662968099861
    $lvl1_prop_data[3,9] = 0.18294701986754966
    $lvl1_prop_data[4,0] = "GC_Cumulus"
    $lvl1_prop_data[4,1] = 0.41935483870967744
    $lvl1_prop_data[4,2] = 0.10964912280701754
    $lvl1_prop_data[4,3] = 0.34074074074074073
    $lvl1_prop_data[4,4] = 0.2088888888888889
    $lvl1_prop_data[4,5] = 0.19377162629757785
    $lvl1_prop_data[4,6] = 0.17391304347826086
    $lvl1_prop_data[4,7] = 0.25256673511293637
    $lvl1_prop_data[4,8] = 0.1927877947295423
    $lvl1_prop_data[4,9] = 0.21688741721854304
    $lvl1_prop_data[5,0] = "GC_Luteinizing and CL"
    $lvl1_prop_data[5,1] = 2.846774193548387
    $lvl1_prop_data[5,2] = 0.9605263157894737
    $lvl1_prop_data[5,3] = 1.5333333333333334
    $lvl1_prop_data[5,4] = 1.3111111111111111
    $lvl1_prop_data[5,5] = 0.5051903114186851
    $lvl1_prop_data[5,6] = 0.7584541062801933
    $lvl1_prop_data[5,7] = 1.5995893223819302
    $lvl1_prop_data[5,8] = 0.8294036061026352
    $lvl1_prop_data[5,9] = 1.1399006622516556
    $lvl1_prop_data[6,0] = "GC_Mitotic"
    $lvl1_prop_data[6,1] = 1.0
    $lvl1_prop_data[6,2] = 1.0
    $lvl1_prop_data[6,3] = 1.0
    $lvl1_prop_data[6,4] = 1.0
    $lvl1_prop_data[6,5] = 1.0
    $lvl1_prop_data[6,6] = 1.0
    $lvl1_prop_data[6,7] = 1.0
    $lvl1_prop_data[6,8] = 1.0
    $lvl1_prop_data[6,9] = 1.0
    $lvl1_prop_data[7,0] = "GC_Mural"
    $lvl1_prop_data[7,1] = 3.314516129032258
    $lvl1_prop_data[7,2] = 2.793859649122807
    $lvl1_prop_data[7,3] = 3.6074074074074076
    $lvl1_prop_data[7,4] = 1.8977777777777778
    $lvl1_prop_data[7,5] = 1.9169550173010381
    $lvl1_prop_data[7,6] = 2.449275362318841
    $lvl1_prop_data[7,7] = 3.1519507186858315
    $lvl1_prop_data[7,8] = 2.0638002773925104
    $lvl1_prop_data[7,9] = 2.502483443708609
    $lvl1_prop_data[8,0] = "Granulosa"
    $lvl1_prop_data[8,1] = 0.0
    $lvl1_prop_data[8,2] = 0.0
    $lvl1_prop_data[8,3] = 0.0
    $lvl1_prop_data[8,4] = 0.0
    $lvl1_prop_data[8,5] = 0.0
    $lvl1_prop_data[8,6] = 0.0
    $lvl1_prop_data[8,7] = 0.0
    $lvl1_prop_data[8,8] = 0.0
    $lvl1_prop_data[8,9] = 0.0
    $lvl1_prop_data[9,0] = "Immune"
    $lvl1_prop_data[9,1] = 0.18548387096774194
    $lvl1_prop_data[9,2] = 0.06140350877192982
    $lvl1_prop_data[9,3] = 0.08888888888888889
    $lvl1_prop_data[9,4] = 0.07555555555555556
    $lvl1_prop_data[9,5] = 0.04498269896193772
    $lvl1_prop_data[9,6] = 0.0821256038647343
    $lvl1_prop_data[9,7] = 0.10061601642710473
    $lvl1_prop_data[9,8] = 0.0651872399445215
    $lvl1_prop_data[9,9] = 0.07947019867549669
    $lvl1_prop_data[10,0] = "M_Early Theca"
    $lvl1_prop_data[10,1] = 0.8145161290322581
    $lvl1_prop_data[10,2] = 0.5657894736842105
    $lvl1_prop_data[10,3] = 0.7555555555555555
    $lvl1_prop_data[10,4] = 0.4888888888888889
    $lvl1_prop_data[10,5] = 0.7093425605536332
    $lvl1_prop_data[10,6] = 0.4830917874396135
    $lvl1_prop_data[10,7] = 0.6817248459958932
    $lvl1_prop_data[10,8] = 0.5755894590846047
    $lvl1_prop_data[10,9] = 0.6183774834437086
    $lvl1_prop_data[11,0] = "M_Fibroblast-like Stroma"
    $lvl1_prop_data[11,1] = 0.12096774193548387
    $lvl1_prop_data[11,2] = 0.039473684210526314
    $lvl1_prop_data[11,3] = 0.05185185185185185
    $lvl1_prop_data[11,4] = 0.057777777777777775
    $lvl1_prop_data[11,5] = 0.031141868512110725
    $lvl1_prop_data[11,6] = 0.028985507246376812
    $lvl1_prop_data[11,7] = 0.06365503080082136
    $lvl1_prop_data[11,8] = 0.038834951456310676
    $lvl1_prop_data[11,9] = 0.048841059602649006
    $lvl1_prop_data[12,0] = "M_Pericyte"
    $lvl1_prop_data[12,1] = 0.25806451612903225
    $lvl1_prop_data[12,2] = 0.12719298245614036
    $lvl1_prop_data[12,3] = 0.3925925925925926
    $lvl1_prop_data[12,4] = 0.07111111111111111
    $lvl1_prop_data[12,5] = 0.09342560553633218
    $lvl1_prop_data[12,6] = 0.13526570048309178
    $lvl1_prop_data[12,7] = 0.23408624229979466
    $lvl1_prop_data[12,8] = 0.09847434119278779
    $lvl1_prop_data[12,9] = 0.15314569536423842
    $lvl1_prop_data[13,0] = "M_Smooth Muscle"
    $lvl1_prop_data[13,1] = 0.16129032258064516
    $lvl1_prop_data[13,2] = 0.06140350877192982
    $lvl1_prop_data[13,3] = 0.14814814814814814
    $lvl1_prop_data[13,4] = 0.06222222222222222
    $lvl1_prop_data[13,5] = 0.06920415224913495
    $lvl1_prop_data[13,6] = 0.0821256038647343
    $lvl1_prop_data[13,7] = 0.11088295687885011
    $lvl1_prop_data[13,8] = 0.07073509015256588
    $lvl1_prop_data[13,9] = 0.0869205298013245
    $lvl1_prop_data[14,0] = "M_Steroidogenic Stroma"
    $lvl1_prop_data[14,1] = 3.0403225806451615
    $lvl1_prop_data[14,2] = 0.9429824561403509
    $lvl1_prop_data[14,3] = 2.214814814814815
    $lvl1_prop_data[14,4] = 1.0266666666666666
    $lvl1_prop_data[14,5] = 1.1038062283737025
    $lvl1_prop_data[14,6] = 1.0144927536231885
    $lvl1_prop_data[14,7] = 1.8295687885010268
    $lvl1_prop_data[14,8] = 1.0540915395284327
    $lvl1_prop_data[14,9] = 1.3667218543046358
    $lvl1_prop_data[15,0] = "M_Steroidogenic Theca"
    $lvl1_prop_data[15,1] = 1.814516129032258
    $lvl1_prop_data[15,2] = 0.6754385964912281
    $lvl1_prop_data[15,3] = 1.4222222222222223
    $lvl1_prop_data[15,4] = 1.0133333333333334
    $lvl1_prop_data[15,5] = 0.7058823529411765
    $lvl1_prop_data[15,6] = 0.9130434782608695
    $lvl1_prop_data[15,7] = 1.1724845995893223
    $lvl1_prop_data[15,8] = 0.8613037447988904
    $lvl1_prop_data[15,9] = 0.9867549668874173
    $lvl1_prop_data[16,0] = "Mesenchyme"
    $lvl1_prop_data[16,1] = 0.0
    $lvl1_prop_data[16,2] = 0.0
    $lvl1_prop_data[16,3] = 0.0
    $lvl1_prop_data[16,4] = 0.0
    $lvl1_prop_data[16,5] = 0.0
    $lvl1_prop_data[16,6] = 0.0
    $lvl1_prop_data[16,7] = 0.0
    $lvl1_prop_data[16,8] = 0.0
    $lvl1_prop_data[16,9] = 0.0
    $lvl1_prop_data[17,0] = "Total"
    $lvl1_prop_data[17,1] = 14.60483870967742
    $lvl1_prop_data[17,2] = 7.62280701754386
    $lvl1_prop_data[17,3] = 12.192592592592593
    $lvl1_prop_data[17,4] = 7.44
    $lvl1_prop_data[17,5] = 6.58477508650519
    $lvl1_prop_data[17,6] = 7.555555555555555
    $lvl1_prop_data[17,7] = 10.66735112936345
    $lvl1_prop_data[17,8] = 7.130374479889043
    $lvl1_prop_data[17,9] = 8.556291390728477
$level1prop.Range("A1:J18").Value = $lvl1_prop_data

$level1prop.Range("A1:J1").Font.Bold = $true
$level1prop.Range("A1:J1").HorizontalAlignment = -4108

$level0.Activate()
$level0.Range("A1").Select()

Write-Output "edit complete"
